# Adds yuv_* color statistic columns (BG:BX) to the sheet, for the "hsv lab xyz yuv" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels for columns BG..BX (row 1) ---
$headers = @("yuv_y_mean","yuv_y_median","yuv_y_variance","yuv_y_std_dev","yuv_y_25","yuv_y_75","yuv_u_mean","yuv_u_median","yuv_u_variance","yuv_u_std_dev","yuv_u_25","yuv_u_75","yuv_v_mean","yuv_v_median","yuv_v_variance","yuv_v_std_dev","yuv_v_25","yuv_v_75")

$nHeaderCols = $headers.Count
$headerArr = New-Object 'object[,]' 1,$nHeaderCols
for ($c = 0; $c -lt $nHeaderCols; $c++) {
    $headerArr[0,$c] = $headers[$c]
}
$ws.Range("BG1:BX1").Value = $headerArr

# Match the existing header formatting (bold, centered, thin border) by
# copying the format of the last existing header cell (BF1).
$ws.Range("BF1").Copy() | Out-Null
$ws.Range("BG1:BX1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- New numeric data for columns BG..BX, rows 2..19 ---
$dataRows = @(
    @(69.97129509581187,63,1205.85327851328,34.72539817645408,47,85,119.1857191406426,128,130.6443139480517,11.42997436340308,111,128,141.4154584426911,128,275.5525333936407,16.59977510069461,128,155),
    @(65.32000986781715,54,1554.67485639792,39.42936540698975,38,86,126.088028271993,128,36.09171766507369,6.007638276816747,128,128,130.9728139241537,128,82.59342989819632,9.088092753608777,128,128),
    @(69.82749431715365,62,1131.780086437595,33.64193939768626,47,86,119.1706267009532,128,140.8683125413667,11.86879574941648,109,128,139.8080439194962,128,235.6229570217146,15.35001488669358,128,153),
    @(71.54639323103399,66,1272.606619256486,35.67361236623628,47,90,118.0092790467851,128,191.6864400380452,13.84508721670056,106,128,142.4714237472749,128,382.7622217372686,19.56430989678064,128,159),
    @(72.82136670404169,66,1086.084923268611,32.95580257357739,48,90,125.6385803222656,128,56.4764988735939,7.515084754917532,128,128,131.1469297409058,128,93.63746895053286,9.676645542259614,128,128),
    @(88.82008592357541,90,826.5811705849128,28.75032470399096,65,110,116.8296035844991,118,141.7173508056759,11.90450968354749,107,128,141.7083462735575,143,200.004844729019,14.14230690973078,128,154),
    @(127.9930184696415,121,1474.238022030635,38.39580734963956,94,164,116.370115283845,128,224.3214711505171,14.97736529402008,104,128,139.649122155822,128,214.4875230554002,14.64539255381706,128,153),
    @(117.4520811956279,117,780.2514096399237,27.932980679475,96,136,125.4886566797892,128,64.01310769281541,8.000819188859063,128,128,130.5319089889526,128,58.84541545595001,7.671076551302953,128,128),
    @(135.0340774868683,132,1133.373841594867,33.66561809316543,106,162,124.3918892542521,128,105.6732123504673,10.27974767931914,128,128,131.5215199788412,128,97.3456570890782,9.866390276543807,128,128),
    @(137.0639757846268,137,888.5537392587412,29.80861853992468,114,160,116.2826135105205,128,219.7868331843283,14.82520938079218,104,128,139.2445594457222,128,185.1543159579124,13.60714209369155,128,152),
    @(156.5852682654553,159,807.3982340355473,28.41475380916659,133,178,116.5739208495892,128,201.1173043922009,14.18158328227849,107,128,139.287203581003,128,151.6190122722358,12.31336721909307,128,151),
    @(149.8006256499033,149,681.7126418030423,26.10962737771342,130,168,125.3196557362874,128,79.29437642872357,8.904738987119362,128,128,130.5768473943075,128,59.85146358300949,7.736372766549547,128,128),
    @(157.7213024228718,153,627.2485987559955,25.04493159814967,140,177,124.4005651473999,128,77.60302075553871,8.809257673353567,128,128,131.3976256052653,128,64.62705076382015,8.039095145837008,128,128),
    @(154.178888702525,154,669.1862034419659,25.86863358281542,134,172,125.8896697362264,128,45.33790740511649,6.733342959118931,128,128,130.1271864573161,128,41.322547741428,6.428261642265971,128,128),
    @(151.3998498170306,152,670.68821322316,25.89764879720088,132,172,126.0448191960653,128,46.60446166558746,6.826746052519272,128,128,129.6372483571371,128,30.79283573493937,5.549129277187491,128,128),
    @(154.9050826670551,159,783.8468467039286,27.99726498613621,139,176,117.8986150616364,121,144.9803849716564,12.0407800815253,110,128,138.1891388170051,141,111.7593856376652,10.57163117204082,128,146),
    @(157.3742395636669,158,603.0285125064859,24.55663886826709,144,171,126.2380634943644,128,39.12484782414402,6.254985837245679,128,128,129.7914355595907,128,33.31207421048003,5.771661304206964,128,128),
    @(164.8959850843257,170,613.2283322755437,24.76344750384211,152,182,124.3725884755452,128,73.07829577324392,8.548584430959545,128,128,131.0718037287394,128,50.01187780150411,7.071907649390234,128,128)
)

$nRows = $dataRows.Count
$nCols = $dataRows[0].Count
$dataArr = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $dataArr[$r,$c] = $dataRows[$r][$c]
    }
}
$ws.Range("BG2:BX19").Value = $dataArr
